# "Updated few failed test cases of Regression test run"
#
# Updates a handful of expected/current-value cells on the two "Loop B"
# regression sheets, and moves the active sheet/selection from
# "Add Devices Loop A" (G1) over to "Add_Devices_LoopB_FIM" (A9), leaving
# a new selection behind on each of the other sheets as well - matching
# where the tester's cursor ended up on each tab after the run.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Add Devices Loop A")
$ws2 = $wb.Worksheets.Item("Add_Devices_LoopB_PFI")
$ws3 = $wb.Worksheets.Item("Add_Devices_LoopB_FIM")

# --- Fix up failing actual-value cells -------------------------------
$ws2.Range("G2").Value = 427.2
$ws2.Range("G3").Value = 427.2

$ws3.Range("G2").Value = 360.6
$ws3.Range("G3").Value = 344.6
$ws3.Range("J9").Value = 17

# --- Leave each sheet's selection where the tester left it -----------
$ws1.Activate()
$ws1.Range("C25").Select()

$ws2.Activate()
$ws2.Range("G3").Select()

# Loop B FIM ends up the active tab when the workbook is saved
$ws3.Activate()
$ws3.Range("A9").Select()
